# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.193.75"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +5.78%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.782.89"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.36%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "243.99"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4919"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2668"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.51%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06244"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.777.06"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.64%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "16.48"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.15%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07017"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.45%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.6254"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.64%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.627"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.03%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "79.85"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.47%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "28.162.33"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +6.34%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000007216"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.11%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.02"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +5.71%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "2.006.50"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.83%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.553"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.15%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.732"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.39%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.225"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.78%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "141.39"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.63%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "15.78"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.00%  "

$ws.Range("E27").Value = "  +5.50%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "109.10"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.91%  "

$ws.Range("E29").Value = "  +0.34%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.168"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +6.69%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08252"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.63%  "

$ws.Range("E32").Value = "  +2.94%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.04881"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +8.66%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.069"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +6.99%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.611"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.6505"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.36%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.9432"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.30%  "

$ws.Range("E38").Value = "  +7.37%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.041"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.935"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +6.93%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.01546"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.42%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.15%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "99.69"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.48%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.3978"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.46%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "7.162"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +4.64%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.1200"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.00%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.05431"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.93%  "

$ws.Range("E48").Value = "  +2.52%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.293"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +5.45%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "30.58"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.25%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "52.75"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.43%  "
